{"js": "// 2022 Tiffin Allegro Open Road 32SA listing \u2014 price correction, mileage\n// caveat, and suspension-upgrade price correction.\n//\n// Four places in the document's visible text actually change:\n//   1. Headline price: $156,000 -> $154,500\n//   2. Mileage bullet: \"Mileage: 17,250\" -> \"Mileage: ~18,250 (it changes monthly)\"\n//   3. Suspension bullet: \"$28,000 upgrade\" -> \"$19,000 upgrade\"\n//   4. Asking-price bullet: $156,000 -> $154,500\n//\n// (The rest of the canonical-XML diff is just the spell-checker re-running\n// and splitting runs / adding proofErr markers around words like\n// \"LiquidSpring\", \"Valterra\", \"Splendide\" \u2014 the visible text there is\n// unchanged, so nothing to do for those paragraphs.)\n\nasync function replaceEverywhere(body, findText, replaceText) {\n  const hits = body.search(findText, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length === 0) {\n    throw new Error(`Find/Replace failed: could not find '${findText}'`);\n  }\n  for (const hit of hits.items) {\n    hit.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1 & 4: both occurrences of the old headline/asking price -> new price.\nawait replaceEverywhere(body, \"$156,000\", \"$154,500\");\n\n// 2: mileage bullet gains a leading \"~\" on 18 and a trailing caveat.\nawait replaceEverywhere(body, \"Mileage: 17,250\", \"Mileage: ~18,250 (it changes monthly)\");\n\n// 3: suspension upgrade price.\nawait replaceEverywhere(body, \"$28,000 upgrade\", \"$19,000 upgrade\");\n", "ps1": "# 2022 Tiffin Allegro Open Road 32SA listing - price correction, mileage\n# caveat, and suspension-upgrade price correction.\n#\n# Four places in the document's visible text actually change:\n#   1. Headline price: $156,000 -> $154,500\n#   2. Mileage bullet: \"Mileage: 17,250\" -> \"Mileage: ~18,250 (it changes monthly)\"\n#   3. Suspension bullet: \"$28,000 upgrade\" -> \"$19,000 upgrade\"\n#   4. Asking-price bullet: $156,000 -> $154,500\n#\n# (The rest of the canonical-XML diff is just the spell-checker re-running\n# and splitting runs / adding proofErr markers around words like\n# \"LiquidSpring\", \"Valterra\", \"Splendide\" - the visible text there is\n# unchanged, so there's nothing else to do for those paragraphs.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Everywhere([string]$findText, [string]$replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n  if (-not $ok) {\n    throw \"Find/Replace failed: could not find '$findText'\"\n  }\n}\n\n# 1 & 4: both occurrences of the old headline/asking price -> new price.\nReplace-Everywhere '$156,000' '$154,500'\n\n# 2: mileage bullet gains a leading \"~\" on 18 and a trailing caveat.\nReplace-Everywhere 'Mileage: 17,250' 'Mileage: ~18,250 (it changes monthly)'\n\n# 3: suspension upgrade price.\nReplace-Everywhere '$28,000 upgrade' '$19,000 upgrade'\n"}
